$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Insert the four "moved" paragraphs (blue / accent1) right after
#    the "15 mm lip" paragraph that precedes "drift eliminator access
#    hatch" -- i.e. before that paragraph, in their new order:
#       Tank hatch cut out
#       250 mm
#       Legs
#       Move to other side
# ------------------------------------------------------------------

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "15 mm lip" + [char]13 -and $p.Next().Range.Text -eq "drift eliminator access hatch" + [char]13) {
        $target = $p
        break
    }
}

$targetEnd = $target.Range.End
$insertRange = $d.Range($targetEnd - 1, $targetEnd - 1)   # just before the paragraph mark

$newXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:color w:val="4472C4" w:themeColor="accent1"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="4472C4" w:themeColor="accent1"/></w:rPr><w:t xml:space="preserve">Tank hatch cut </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:color w:val="4472C4" w:themeColor="accent1"/></w:rPr><w:t>out</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:rPr><w:color w:val="4472C4" w:themeColor="accent1"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="4472C4" w:themeColor="accent1"/></w:rPr><w:t>250 mm</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:color w:val="4472C4" w:themeColor="accent1"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="4472C4" w:themeColor="accent1"/></w:rPr><w:t>Legs</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:rPr><w:color w:val="4472C4" w:themeColor="accent1"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="4472C4" w:themeColor="accent1"/></w:rPr><w:t xml:space="preserve">Move to other </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:color w:val="4472C4" w:themeColor="accent1"/></w:rPr><w:t>sid</w:t></w:r><w:r><w:rPr><w:color w:val="4472C4" w:themeColor="accent1"/></w:rPr><w:t>e</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>
'@

$insertRange.InsertXML($newXml)

# ------------------------------------------------------------------
# 2. Remove the old (red) copies of those four paragraphs that used
#    to sit right after "Add side panels" / before "Slotted panel".
# ------------------------------------------------------------------

$legsPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Legs" + [char]13 -and $p.Range.Font.TextColor.RGB -eq 255) {
        $legsPara = $p
        break
    }
}

$endPara = $legsPara.Next().Next().Next()   # Move to other side / Tank hatch cut out / 250 mm

$deleteRange = $d.Range($legsPara.Range.Start, $endPara.Range.End)
$deleteRange.Delete()
